$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "54.759.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.27%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.340.94"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.91%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "470.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.71"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.40%  "

$ws.Range("E7").Value = "  +0.46%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.500"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.345.36"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.38%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0957"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.93%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.37"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.27%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.317"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.09%  "

$ws.Range("E13").Value = "  +0.37%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.754.78"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.56%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "55.158.81"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.47%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.90"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.72%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000129"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.11%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.346.75"
$ws.Range("D18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.49"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.86%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "310.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.42%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.48%  "

$ws.Range("E22").Value = "  +0.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.72%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "56.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.97%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.995"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.390"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.54%  "

$ws.Range("E27").Value = "  -8.29%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.455.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.34%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.12"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.54%  "

$ws.Range("E30").Value = "  +0.08%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0746"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.44%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "148.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.77%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.90"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.09%  "

$ws.Range("E34").Value = "  -3.60%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.98"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.90%  "

$ws.Range("E36").Value = "  -6.33%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.51"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.49%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.815"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.40%  "

$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "33.35"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.54%  "

$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.998"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.53%  "

$ws.Range("E41").Value = "  -2.19%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.33"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.61%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0939"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.01%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0523"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.80%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.573"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.80%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.14"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.49%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "251.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.53%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0220"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.13%  "

$ws.Range("E49").Value = "  -9.52%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.69"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.60%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.762.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.40%  "
